$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the existing date-cell formatting (column A) down into the new rows
$ws.Range("A41").Copy() | Out-Null
$ws.Range("A42:A44").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New log entries appended to the effort sheet
$ws.Range("A42").Value = 41228
$ws.Range("B42").Value = 2.5
$ws.Range("D42").Value = "Test case tc08 is working well, a fix of rtos.c was required. All test cases rerun. Manual updated"

$ws.Range("A43").Value = 41229
$ws.Range("B43").Value = 3.25
$ws.Range("D43").Value = "Preparation of release, new test case tc09"

$ws.Range("A44").Value = 41232
$ws.Range("B44").Value = 1.25
$ws.Range("D44").Value = "Makefile: Workaround for 12 Bit Branch distance problem with core.a"

# Update the view so the newly added row is visible/selected (scroll the
# window down so row 38 is the top visible row, then select the new entry)
$excel.ActiveWindow.ScrollRow = 38
$ws.Range("A44").Select() | Out-Null
